# Phieu xuat kho nguyen lieu san xuat - update SanPham.xlsx template row 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: "Loai san pham" was blank -> now filled with SP_NHA_CUNG_CAP
$ws.Range("E2").Value = "SP_NHA_CUNG_CAP"

# I2 / J2: "Ty le chiet khau" / "Muc loi nhuan" were text ("5.00"/"10.00") -> now numeric
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 10

# Move the active selection to K11 (matches the saved view state in the workbook)
[void]$ws.Range("K11").Select()
